# "move BC temporary files to repository files"
#
# The source temp file had a provisional 0.99 "almost fully phased-in"
# market-share-class-min value for the 2035-2050 columns (T3:W3); moving it
# into the repo finalizes those years to the full 1 (100%).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core data edit: ZEV market share_class_min reaches 100% (1) for 2035
# onward instead of the provisional 0.99.
$ws.Range("T3:W3").Value = 1

# Leave the sheet selected over its full used range, matching the saved
# view state of the finalized file.
$ws.Range("A1:X8").Select()
